$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update match rows 49-52 and 91-93: odds/teams data was corrected/reordered. ---
# Columns A (index) and E (date) are unchanged; only F:V are updated.
# Row 49
$ws.Range('F49').Value = 'Pogon Siedlce'
$ws.Range('G49').Value = 3
$ws.Range('H49').Value = 'Sandecja Nowy S.'
$ws.Range('I49').Value = 0
$ws.Range('J49').Value = 2.28
$ws.Range('K49').Value = '26/08/2023 13:13'
$ws.Range('L49').Value = 2.77
$ws.Range('M49').Value = '26/08/2023 16:56'
$ws.Range('N49').Value = 3.08
$ws.Range('O49').Value = '26/08/2023 13:13'
$ws.Range('P49').Value = 3.12
$ws.Range('Q49').Value = '26/08/2023 16:55'
$ws.Range('R49').Value = 3.09
$ws.Range('S49').Value = '26/08/2023 13:13'
$ws.Range('T49').Value = 2.56
$ws.Range('U49').Value = '26/08/2023 16:56'
$ws.Range('V49').Value = 'https://www.betexplorer.com/football/poland/division-2/pogon-siedlce-sandecja-nowy-s/bT3PVrIQ/'

# Row 50
$ws.Range('F50').Value = 'GKS Jastrzebie'
$ws.Range('G50').Value = 4
$ws.Range('H50').Value = 'Polonia Bytom'
$ws.Range('I50').Value = 2
$ws.Range('J50').Value = 1.85
$ws.Range('K50').Value = '26/08/2023 13:13'
$ws.Range('L50').Value = 1.83
$ws.Range('M50').Value = '26/08/2023 16:58'
$ws.Range('N50').Value = 3.42
$ws.Range('O50').Value = '26/08/2023 13:13'
$ws.Range('P50').Value = 3.6
$ws.Range('Q50').Value = '26/08/2023 16:58'
$ws.Range('R50').Value = 4.09
$ws.Range('S50').Value = '26/08/2023 13:13'
$ws.Range('T50').Value = 4.03
$ws.Range('U50').Value = '26/08/2023 16:58'
$ws.Range('V50').Value = 'https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-polonia-bytom/464LWO2K/'

# Row 51
$ws.Range('F51').Value = 'Zaglebie II'
$ws.Range('G51').Value = 2
$ws.Range('H51').Value = 'Olimpia Elblag'
$ws.Range('I51').Value = 1
$ws.Range('J51').Value = 2.01
$ws.Range('K51').Value = '26/08/2023 13:13'
$ws.Range('L51').Value = 2.44
$ws.Range('M51').Value = '26/08/2023 16:30'
$ws.Range('N51').Value = 3.27
$ws.Range('O51').Value = '26/08/2023 13:13'
$ws.Range('P51').Value = 3.34
$ws.Range('Q51').Value = '26/08/2023 16:30'
$ws.Range('R51').Value = 3.5
$ws.Range('S51').Value = '26/08/2023 13:13'
$ws.Range('T51').Value = 2.77
$ws.Range('U51').Value = '26/08/2023 16:30'
$ws.Range('V51').Value = 'https://www.betexplorer.com/football/poland/division-2/zaglebie-olimpia-elblag/QF5HX4mE/'

# Row 52
$ws.Range('F52').Value = 'Wisla Pulawy'
$ws.Range('G52').Value = 0
$ws.Range('H52').Value = 'Hutnik Krakow'
$ws.Range('I52').Value = 1
$ws.Range('J52').Value = 1.78
$ws.Range('K52').Value = '26/08/2023 13:13'
$ws.Range('L52').Value = 2.04
$ws.Range('M52').Value = '26/08/2023 16:56'
$ws.Range('N52').Value = 3.61
$ws.Range('O52').Value = '26/08/2023 13:13'
$ws.Range('P52').Value = 3.28
$ws.Range('Q52').Value = '26/08/2023 16:54'
$ws.Range('R52').Value = 3.99
$ws.Range('S52').Value = '26/08/2023 13:13'
$ws.Range('T52').Value = 3.62
$ws.Range('U52').Value = '26/08/2023 16:56'
$ws.Range('V52').Value = 'https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-hutnik-krakow/YLfiR2Qs/'

# Row 91
$ws.Range('F91').Value = 'Sandecja Nowy S.'
$ws.Range('G91').Value = 0
$ws.Range('H91').Value = 'KKS Kalisz'
$ws.Range('I91').Value = 2
$ws.Range('J91').Value = 2.75
$ws.Range('K91').Value = '29/09/2023 02:12'
$ws.Range('L91').Value = 2.78
$ws.Range('M91').Value = '30/09/2023 14:41'
$ws.Range('N91').Value = 3.09
$ws.Range('O91').Value = '29/09/2023 02:12'
$ws.Range('P91').Value = 3.12
$ws.Range('Q91').Value = '30/09/2023 14:41'
$ws.Range('R91').Value = 2.35
$ws.Range('S91').Value = '29/09/2023 02:12'
$ws.Range('T91').Value = 2.56
$ws.Range('U91').Value = '30/09/2023 14:41'
$ws.Range('V91').Value = 'https://www.betexplorer.com/football/poland/division-2/sandecja-nowy-s-kks-kalisz/KbAXndAF/'

# Row 92
$ws.Range('F92').Value = 'Lech Poznan II'
$ws.Range('G92').Value = 0
$ws.Range('H92').Value = 'LKS Lodz II'
$ws.Range('I92').Value = 3
$ws.Range('J92').Value = 2.36
$ws.Range('K92').Value = '29/09/2023 02:12'
$ws.Range('L92').Value = 2.4
$ws.Range('M92').Value = '30/09/2023 14:59'
$ws.Range('N92').Value = 3.25
$ws.Range('O92').Value = '29/09/2023 02:12'
$ws.Range('P92').Value = 3.65
$ws.Range('Q92').Value = '30/09/2023 14:59'
$ws.Range('R92').Value = 2.6
$ws.Range('S92').Value = '29/09/2023 02:12'
$ws.Range('T92').Value = 2.63
$ws.Range('U92').Value = '30/09/2023 14:58'
$ws.Range('V92').Value = 'https://www.betexplorer.com/football/poland/division-2/lech-poznan-lks-lodz/CdkCuE2k/'

# Row 93
$ws.Range('F93').Value = 'Hutnik Krakow'
$ws.Range('G93').Value = 0
$ws.Range('H93').Value = 'Chojniczanka'
$ws.Range('I93').Value = 2
$ws.Range('J93').Value = 2.06
$ws.Range('K93').Value = '29/09/2023 02:12'
$ws.Range('L93').Value = 2.24
$ws.Range('M93').Value = '30/09/2023 14:43'
$ws.Range('N93').Value = 3.21
$ws.Range('O93').Value = '29/09/2023 02:12'
$ws.Range('P93').Value = 3.49
$ws.Range('Q93').Value = '30/09/2023 14:41'
$ws.Range('R93').Value = 3.14
$ws.Range('S93').Value = '29/09/2023 02:12'
$ws.Range('T93').Value = 2.95
$ws.Range('U93').Value = '30/09/2023 14:43'
$ws.Range('V93').Value = 'https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-chojniczanka/t8BTmxe9/'

# --- Append two new match rows (140 and 141) ---
# Copy formatting (bold/border style on col A, date format on col E) from the last existing row.
$ws.Range('A139:V139').Copy($ws.Range('A140:V140'))
$ws.Range('A139:V139').Copy($ws.Range('A141:V141'))

# Row 140
$ws.Range('A140').Value = 139
$ws.Range('B140').Value = 'poland'
$ws.Range('C140').Value = 'division-2'
$ws.Range('D140').Value = '2023-2024'
$ws.Range('E140').Value = 45235.52083333334
$ws.Range('F140').Value = 'Hutnik Krakow'
$ws.Range('G140').Value = 3
$ws.Range('H140').Value = 'Kotwica Kolobrzeg'
$ws.Range('I140').Value = 5
$ws.Range('J140').Value = 2.46
$ws.Range('K140').Value = '04/11/2023 00:42'
$ws.Range('L140').Value = 2.63
$ws.Range('M140').Value = '05/11/2023 12:10'
$ws.Range('N140').Value = 3.14
$ws.Range('O140').Value = '04/11/2023 00:42'
$ws.Range('P140').Value = 3.18
$ws.Range('Q140').Value = '05/11/2023 12:10'
$ws.Range('R140').Value = 2.57
$ws.Range('S140').Value = '04/11/2023 00:42'
$ws.Range('T140').Value = 2.66
$ws.Range('U140').Value = '05/11/2023 12:10'
$ws.Range('V140').Value = 'https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-kotwica-kolobrzeg/t4LHWwLA/'

# Row 141
$ws.Range('A141').Value = 140
$ws.Range('B141').Value = 'poland'
$ws.Range('C141').Value = 'division-2'
$ws.Range('D141').Value = '2023-2024'
$ws.Range('E141').Value = 45235.625
$ws.Range('F141').Value = 'S. Wola'
$ws.Range('G141').Value = 1
$ws.Range('H141').Value = 'Stezyca'
$ws.Range('I141').Value = 2
$ws.Range('J141').Value = 2.03
$ws.Range('K141').Value = '04/11/2023 03:13'
$ws.Range('L141').Value = 2.1
$ws.Range('M141').Value = '05/11/2023 14:57'
$ws.Range('N141').Value = 3.21
$ws.Range('O141').Value = '04/11/2023 03:13'
$ws.Range('P141').Value = 3.17
$ws.Range('Q141').Value = '05/11/2023 14:57'
$ws.Range('R141').Value = 3.22
$ws.Range('S141').Value = '04/11/2023 03:13'
$ws.Range('T141').Value = 3.57
$ws.Range('U141').Value = '05/11/2023 14:54'
$ws.Range('V141').Value = 'https://www.betexplorer.com/football/poland/division-2/stal-stalowa-wola-stezyca/hrT4Zuyh/'
